$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 10.67
$ws.Range("E3").Value = 10.39
$ws.Range("F3").Value = 10.17

# Row 4
$ws.Range("C4").Value = 9.33
$ws.Range("E4").Value = 9.9
$ws.Range("F4").Value = 10.14

# Row 5
$ws.Range("C5").Value = 9.609999999999999
$ws.Range("D5").Value = 10.1
$ws.Range("F5").Value = 10.29
$ws.Range("G5").Value = 9.18

# Row 6
$ws.Range("C6").Value = 9.83
$ws.Range("D6").Value = 9.859999999999999
$ws.Range("E6").Value = 9.710000000000001
$ws.Range("G6").Value = 10.34
$ws.Range("H6").Value = 11.26

# Row 7
$ws.Range("E7").Value = 10.82
$ws.Range("F7").Value = 9.66

# Row 8
$ws.Range("F8").Value = 8.74
$ws.Range("J8").Value = 14.75

# Row 10
$ws.Range("H10").Value = 5.25
